# Apply corrections to the District column (G) so that the various
# misspellings / abbreviations of "Madhugiri" are replaced with the
# official name "Tumakuru (Tumkur)", correct a couple of other district
# names, and drop a handful of stray empty cells in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-25 and 27-38 (i.e. all except row 26): District -> "Tumakuru (Tumkur)"
$tumakuruRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,27,28,29,30,31,32,33,34,35,36,37,38)
foreach ($r in $tumakuruRows) {
    $ws.Range("G$r").Value = "Tumakuru (Tumkur)"
}

# Row 40: "Bangalore North Range -3" -> "Bengaluru (Bangalore) Rural"
$ws.Range("G40").Value = "Bengaluru (Bangalore) Rural"

# Row 49: "Bangalor" -> "Bengaluru (Bangalore) Urban"
$ws.Range("G49").Value = "Bengaluru (Bangalore) Urban"

# Remove stray empty cells in column F for rows 22, 26, 41, 50, 51
$emptyFRows = @(22,26,41,50,51)
foreach ($r in $emptyFRows) {
    $ws.Range("F$r").ClearContents()
}
